# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

# Step 1: rename existing summary sheet '总计' to '2022-Q1'
$total = $wb.Worksheets.Item('总计')
$total.Name = '2022-Q1'

# Step 2: duplicate it (still holding the old 总计 data) right after itself,
# then rename the duplicate back to '总计'. This keeps sheetId/rId allocation
# in line with the expected layout (2022-Q1 keeps rId6, new 总计 gets rId7).
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item($total.Index + 1)
$newTotal.Name = '总计'

# Step 3: clear the renamed '2022-Q1' sheet (old 总计 data) and rebuild it
# as a per-fund holdings table.
$q1 = $total
$q1.Cells.Clear()

$headers6 = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($j = 0; $j -lt $headers6.Count; $j++) {
    $cell = $q1.Cells.Item(1, $j + 2)
    $cell.Value = $headers6[$j]
}
# copy header style (bold/centered/bordered) from an existing sheet
$styleSrc = $wb.Worksheets.Item('2021-Q4')
$styleSrc.Cells.Item(1, 2).Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$data6 = @(
    @('005609', '富国军工主题混合A', '74.07', '92.86', '7.82', '5.7923', 4),
    @('000404', '易方达新兴成长灵活配置混合', '51.67', '91.13', '5.64', '2.9142', 5),
    @('110009', '易方达价值精选混合', '41.58', '88.97', '3.77', '1.5676', 9),
    @('000986', '太平灵活配置混合型发起式', '18.13', '82.45', '4.20', '0.7615', 7),
    @('001877', '宝盈国家安全战略沪港深股票', '13.59', '90.66', '3.73', '0.5069', 7),
    @('009693', '富国积极成长一年定期开放混合', '17.82', '98.74', '2.82', '0.5025', 8),
    @('011113', '富国军工主题混合C', '6.15', '92.86', '7.82', '0.4809', 4),
    @('240001', '华宝宝康消费品混合', '11.24', '62.57', '3.45', '0.3878', 5),
    @('005977', '中信保诚至兴灵活配置混合A', '7.64', '89.13', '4.99', '0.3812', 6),
    @('001825', '建信中国制造2025股票', '8.09', '94.18', '4.54', '0.3673', 10),
    @('001268', '富国国家安全主题混合', '4.43', '90.49', '5.58', '0.2472', 6),
    @('550009', '信诚中小盘混合', '5.23', '87.68', '4.16', '0.2176', 6),
    @('000124', '华宝服务优选混合', '6.61', '89.07', '2.53', '0.1672', 10),
    @('010114', '华宝新兴成长混合', '3.26', '83.47', '4.75', '0.1548', 4),
    @('011506', '建信高端装备股票型证券投资基金A', '2.59', '85.91', '3.90', '0.1010', 8),
    @('005978', '中信保诚至兴灵活配置混合C', '1.83', '89.13', '4.99', '0.0913', 6),
    @('001105', '信达澳银转型创新股票', '2.54', '90.59', '3.50', '0.0889', 10),
    @('004183', '富国产业升级混合', '1.47', '92.47', '5.76', '0.0847', 4),
    @('005876', '易方达鑫转增利混合A', '12.58', '26.96', '0.65', '0.0818', 10),
    @('005674', '诺德消费升级灵活配置混合', '1.38', '93.05', '5.06', '0.0698', 10),
    @('001103', '前海开源工业革命4.0灵活配置混合', '4.24', '37.48', '1.45', '0.0615', 9),
    @('003842', '中邮景泰灵活配置混合A', '5.76', '33.95', '1.01', '0.0582', 10),
    @('011507', '建信高端装备股票型证券投资基金C', '0.90', '85.91', '3.90', '0.0351', 8),
    @('011073', '鹏华安润混合A', '3.07', '29.79', '1.14', '0.0350', 4),
    @('163818', '中银中小盘成长混合', '0.98', '87.49', '2.79', '0.0273', 2),
    @('005877', '易方达鑫转增利混合C', '4.17', '26.96', '0.65', '0.0271', 10),
    @('571002', '诺德主题灵活配置混合', '0.67', '77.03', '3.91', '0.0262', 10),
    @('002152', '华宝核心优势灵活配置混合', '0.45', '90.91', '4.65', '0.0209', 2),
    @('003843', '中邮景泰灵活配置混合C', '0.37', '33.95', '1.01', '0.0037', 10),
    @('011074', '鹏华安润混合C', '0.25', '29.79', '1.14', '0.0028', 4)
)

for ($i = 0; $i -lt $data6.Count; $i++) {
    $row = $data6[$i]
    $r = $i + 2
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}
# copy column-A style (bold/centered/bordered) down the new rows
$styleSrc.Cells.Item(2, 1).Copy()
$q1.Range("A2:A31").PasteSpecial(-4122)

# Step 4: rebuild the new '总计' sheet with the prepended 2022-Q1 row
$tot = $newTotal
$tot.Cells.Clear()
$headers7 = @('日期', '持有数量(只)', '持有市值(亿元)')
for ($j = 0; $j -lt $headers7.Count; $j++) {
    $tot.Cells.Item(1, $j + 2).Value = $headers7[$j]
}
$styleSrc.Cells.Item(1, 2).Copy()
$tot.Range("B1:D1").PasteSpecial(-4122)

$data7 = @(
    @('2022-Q1', 30, 15.26),
    @('2021-Q4', 73, 33.83),
    @('2021-Q3', 39, 21.75),
    @('2021-Q2', 36, 24.62),
    @('2021-Q1', 42, 23.13),
    @('2020-Q4', 23, 10.82)
)
for ($i = 0; $i -lt $data7.Count; $i++) {
    $row = $data7[$i]
    $r = $i + 2
    $tot.Cells.Item($r, 1).Value = $i
    $tot.Cells.Item($r, 2).Value = $row[0]
    $tot.Cells.Item($r, 3).Value = $row[1]
    $tot.Cells.Item($r, 4).Value = $row[2]
}
$styleSrc.Cells.Item(2, 1).Copy()
$tot.Range("A2:A7").PasteSpecial(-4122)

# restore the originally active sheet/tab
$wb.Worksheets.Item('2020-Q4').Activate()

Write-Output "done"